$d = $word.ActiveDocument

# Locate the misspelled word "BackgroundNoiseSupression" (missing a "p")
$word_range = $d.Content
$found = $word_range.Find.Execute("BackgroundNoiseSupression")
if (-not $found) {
    throw "Could not find 'BackgroundNoiseSupression' in the document"
}

$wordStart = $word_range.Start
$wordEnd = $word_range.End

# The run splits as "BackgroundNoiseSu" + "p" + "pression" -- the inserted
# "p" lands right after "BackgroundNoiseSu" (offset 17 into the found word).
$splitOffset = $wordStart + 17

# Insert the missing "p" in place, correcting the typo to "BackgroundNoiseSuppression".
$insertionPoint = $d.Range($splitOffset, $splitOffset)
$insertionPoint.InsertAfter("p")

# The inserted "p" now occupies a single character range; force it (and the
# text that follows it) onto their own runs -- matching how this edit was
# actually saved -- by temporarily bookmarking the new character and then
# removing the bookmark. This produces a genuine run boundary without
# leaving any residual character formatting behind.
$pRange = $d.Range($splitOffset, $splitOffset + 1)
$bookmarkName = "TempSplitMarker"
$d.Bookmarks.Add($bookmarkName, $pRange)
$d.Bookmarks($bookmarkName).Delete()
